$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell: "Producto 1" -> "Producto 12"
$ws.Range("B2").Value = "Producto 12"

# Row 2 numeric updates
$ws.Range("C2").Value = 35.0
$ws.Range("E2").Value = 627.0

# Row 3 numeric updates
$ws.Range("C3").Value = 90.0
$ws.Range("E3").Value = 2255.07

# Row 4 numeric updates
$ws.Range("C4").Value = 408.0
$ws.Range("E4").Value = 18.81

# Row 5 numeric updates
$ws.Range("C5").Value = 948.0
$ws.Range("E5").Value = 1879.22

# Row 6 numeric update
$ws.Range("E6").Value = 110.0

# New row 7
$ws.Range("A7").Value = 75.0
$ws.Range("B7").Value = "prueba 54"
$ws.Range("C7").Value = 567.0
$ws.Range("D7").Value = "promo 2x1"
$ws.Range("E7").Value = 800.0
